$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.285.09"
$ws.Range("E2").Value = "  +2.35%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.350.78"
$ws.Range("E3").Value = "  +6.36%  "

$ws.Range("E4").Value = "  -0.20%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "110.32"
$ws.Range("E5").Value = "  +2.37%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "312.91"
$ws.Range("E6").Value = "  +5.47%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.642"
$ws.Range("E7").Value = "  +2.96%  "

$ws.Range("E8").Value = "  -0.18%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.635"
$ws.Range("E9").Value = "  +6.66%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.90"
$ws.Range("E10").Value = "  -1.50%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0940"
$ws.Range("E11").Value = "  +3.43%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.95"
$ws.Range("E12").Value = "  +2.57%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.06"
$ws.Range("E13").Value = "  +11.73%  "

$ws.Range("E14").Value = "  +2.38%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.26"
$ws.Range("E15").Value = "  +9.31%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.707.03"
$ws.Range("E16").Value = "  +6.44%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.347.15"
$ws.Range("E17").Value = "  +5.46%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.275.48"
$ws.Range("E18").Value = "  +2.72%  "

$ws.Range("E19").Value = "  +4.28%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.31"
$ws.Range("E20").Value = "  +0.41%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "75.47"
$ws.Range("E21").Value = "  +4.42%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.60"
$ws.Range("E22").Value = "  +14.60%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.45"
$ws.Range("E23").Value = "  -0.98%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "254.01"
$ws.Range("E24").Value = "  +11.84%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.09"
$ws.Range("E25").Value = "  +1.10%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.06"
$ws.Range("E26").Value = "  +4.58%  "

$ws.Range("E27").Value = "  -1.05%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "39.57"
$ws.Range("E28").Value = "  +3.95%  "

$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.50"
$ws.Range("E29").Value = "  +8.17%  "

$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "174.84"
$ws.Range("E30").Value = "  +0.98%  "

$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.16"
$ws.Range("E31").Value = "  -3.36%  "

$ws.Range("E32").Value = "  -0.81%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0934"
$ws.Range("E33").Value = "  +6.53%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.99"
$ws.Range("E34").Value = "  +8.81%  "

$ws.Range("E35").Value = "  +6.64%  "

$ws.Range("E36").Value = "  -0.88%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0379"
$ws.Range("E37").Value = "  +5.30%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.15"
$ws.Range("E38").Value = "  -2.53%  "

$ws.Range("E39").Value = "  +2.43%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.70"
$ws.Range("E40").Value = "  +11.36%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "72.91"
$ws.Range("E41").Value = "  +3.21%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.48"
$ws.Range("E42").Value = "  +14.85%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.234"
$ws.Range("E43").Value = "  +2.53%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.00"
$ws.Range("E44").Value = "  +3.32%  "

$ws.Range("E45").Value = "  -0.02%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.65"
$ws.Range("E46").Value = "  +4.63%  "

$ws.Range("E47").Value = "  +10.77%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "110.97"
$ws.Range("E48").Value = "  +8.43%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.31"
$ws.Range("E49").Value = "  -0.67%  "

$ws.Range("E50").Value = "  +3.51%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "69.85"
$ws.Range("E51").Value = "  +5.17%  "
